# loginXpaths.xlsx update — "Add files via upload"
#
# Turns the A-only "Combos" xpath list into a 3-column xpath reference sheet:
#   A: LoginFields  (was "Combos")
#   B: LoginButton  (unchanged content, already present)
#   C: RangeButton  (new column with 5 new xpath rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "LoginFields"
$ws.Cells.Item(1, 2).Value = "LoginButton"
$ws.Cells.Item(1, 3).Value = "RangeButton"

# Give the new header cell (C1) the same bold font as B1, but without the
# centered alignment B1 has (matches the new 4th cellXfs entry).
$ws.Cells.Item(1, 3).Font.Bold = $true

# --- New column C data (RangeButton xpaths) --------------------------------
$rangeButtonXpaths = @(
    "/html/body/div[1]/div/div[1]/div/div[3]/div/div/div[1]/div[1]/div[1]/div/div[3]/div[1]/div[2]/div[3]/div[2]/div[1]/div[1]/div/span",
    "/html/body/div[1]/div/div[1]/div/div[3]/div/div/div[1]/div[1]/div[1]/div/div[3]/div[1]/div[2]/div[3]/div[2]/div[1]/div[1]/div",
    "/html/body/div[1]/div/div[1]/div/div[3]/div/div/div[1]/div[1]/div[1]/div/div[3]/div[1]/div[2]/div[3]/div[2]/div[1]/div[2]",
    "/html/body/div[1]/div/div[1]/div/div[3]/div/div/div[1]/div[1]/div[1]/div/div[3]/div[1]/div[2]/div[3]/div[2]/div[1]/div[1]",
    "/html/body/div[1]/div/div[1]/div/div[3]/div/div/div[1]/div[1]/div[1]/div/div[3]/div[1]/div[2]/div[3]/div[2]/div[1]"
)

for ($i = 0; $i -lt $rangeButtonXpaths.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $rangeButtonXpaths[$i]
}

# --- Column sizing (best-fit approximation for the two new text columns) ---
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(3).AutoFit()

# --- Selection / view state -------------------------------------------------
$ws.Range("D1").Select()
